# Insert a new weekly price record at row 511 (pushing the existing
# rows 511-603 down to 512-604), then populate the new row with the
# latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 511:603 down by one row, growing the used range to R604.
$ws.Rows.Item(511).Insert()

# Populate the newly inserted row 511 with the new data point.
$ws.Cells.Item(511, 1).Value = 8
$ws.Cells.Item(511, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(511, 3).Value = 'Coquimbo'
$ws.Cells.Item(511, 4).Value = 45218
$ws.Cells.Item(511, 5).Value = 4
$ws.Cells.Item(511, 6).Value = 100112032
$ws.Cells.Item(511, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(511, 8).Value = 'Sin especificar'
$ws.Cells.Item(511, 9).Value = 'Primera'
$ws.Cells.Item(511, 10).Value = 500
$ws.Cells.Item(511, 11).Value = 13000
$ws.Cells.Item(511, 12).Value = 14000
$ws.Cells.Item(511, 13).Value = 13500
$ws.Cells.Item(511, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(511, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(511, 16).Value = 225
$ws.Cells.Item(511, 17).Value = 60
$ws.Cells.Item(511, 18).Value = 'Hortaliza'
